$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D" = -0.113
    "E" = -0.168
    "G" = 1.102439024390244
    "H" = 1.102439024390244
    "I" = 0.1917073170731708
    "J" = 0.1596703224750368
    "K" = 1.25
    "L" = 0.3048780487804879
    "M" = 3.12
    "N" = 0.1253012048192771
    "O" = 2.496
    "P" = 3.12
    "Q" = 0.1253012048192771
    "R" = 2.496
    "U" = 1.44
    "V" = 0.05783132530120482
    "W" = 0.1201923076923077
    "X" = 0.08648830264224244
    "Y" = 0.03370400505006524
    "Z" = 0.4632768361581921
    "AA" = 0.07397156182459334
    "AB" = 0.08612431860519103
    "AC" = -0.01215275678059768
    "AD" = 0.2
    "AF" = 0.2
    "AG" = -1.24
    "AH" = 0.007968127490039842
    "AI" = 0.02395209580838323
    "AJ" = -0.05240912933220625
    "AK" = -0.1794500723589001
    "AM" = -0.7
    "AN" = 0.2336448598130841
    "AP" = -1.448598130841122
    "AQ" = -1.122857142857143
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
